$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trade data (row 6) matching the existing table layout:
# A=Date, B=Profitable, C=Principle, D=Start Principle, E=BuyPrice,
# F=SellPrice, G=IsShortSell, H=Price Change %, I=Strong trade

# Copy the date-format style (style index 1) from the row above first,
# then set values, so the new cells reuse the existing numFmt (no new
# style gets created in styles.xml).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$ws.Range("A6").Value = 42647.681817129633
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 9971.89
$ws.Range("D6").Value = 10013.450000000001
$ws.Range("E6").Value = 18.12
$ws.Range("F6").Value = 17.97
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = -0.83
$ws.Range("I6").Value = $false
